$wb = $excel.ActiveWorkbook

$wsSetUp   = $wb.Worksheets.Item("SetUp")
$wsLogin   = $wb.Worksheets.Item("LoginPage")
$wsMigr    = $wb.Worksheets.Item("MigrationPage")
$wsHelp    = $wb.Worksheets.Item("HelpAndSupportPage")
$wsHome    = $wb.Worksheets.Item("HomeMovePage")
$wsFeed    = $wb.Worksheets.Item("FeedBackPage")

# ---------------------------------------------------------------------------
# HomeMovePage (sheet5) - main content rewrite
# ---------------------------------------------------------------------------

# Row 2 - addressOnMap now points at "The Dubai Mall"
$wsHome.Range("B2").Value = "The Dubai Mall"

# Row 3 (new) - addressOutside
$wsHome.Range("A3").Value = "addressOutside"
$wsHome.Range("B3").Value = "the world trade center abu dhabi"

# Row 4 (was row 3) - FloorApartmentVilla / 405 (quote-prefixed text, unchanged format)
$wsHome.Range("A4").Value = "FloorApartmentVilla"
$wsHome.Range("B4").Value = "'405"

# Row 5 (new) - TelephoneNumber / 092229774
$wsHome.Range("A5").Value = "TelephoneNumber"
$fA5 = $wsHome.Range("A5").Font
$fA5.Name = "JetBrains Mono"
$fA5.Bold = $true
$fA5.Size = 9.8
$fA5.Color = 14580521

$wsHome.Range("B5").NumberFormat = "@"
$wsHome.Range("B5").HorizontalAlignment = -4131
$wsHome.Range("B5").Value = "092229774"

# Row 6 (new) - TotalAmount / 157.50
$wsHome.Range("A6").Value = "TotalAmount"
$fA6 = $wsHome.Range("A6").Font
$fA6.Name = "JetBrains Mono"
$fA6.Bold = $true
$fA6.Size = 9.8
$fA6.Color = 14580521

$wsHome.Range("B6").NumberFormat = "@"
$wsHome.Range("B6").HorizontalAlignment = -4131
$wsHome.Range("B6").Value = "157.50"

# Row 7 (new) - AccountNumber / Account Number: 042545064
$wsHome.Range("A7").Value = "AccountNumber"
$wsHome.Range("B7").Value = "Account Number: 042545064"
$fB7 = $wsHome.Range("B7").Font
$fB7.Name = "PingFang SC"
$fB7.Size = 12
$wsHome.Rows.Item(7).RowHeight = 19

# Row 8 (new) - Packagename / eLife Unlimited Premium 500
$wsHome.Range("A8").Value = "Packagename"
$fA8 = $wsHome.Range("A8").Font
$fA8.Name = "Calibri Light"
$fA8.Size = 12

$wsHome.Range("B8").Value = "eLife Unlimited Premium 500"
$fB8 = $wsHome.Range("B8").Font
$fB8.Name = "PingFang SC"
$fB8.Size = 12
$wsHome.Rows.Item(8).RowHeight = 19

# Row 9 (new) - InstallationAdress / DUBAI MALL 405 DX
$wsHome.Range("A9").Value = "InstallationAdress"
$wsHome.Range("B9").Value = "DUBAI MALL 405 DX"
$fB9 = $wsHome.Range("B9").Font
$fB9.Name = "PingFang SC"
$fB9.Size = 12
$wsHome.Rows.Item(9).RowHeight = 19

# Row 10 (new) - Moveoutdate / Move-out date: Thu, 25 Mar 2021
$wsHome.Range("A10").Value = "Moveoutdate"
$fA10 = $wsHome.Range("A10").Font
$fA10.Name = "PingFang SC"
$fA10.Size = 12

$wsHome.Range("B10").Value = "Move-out date: Thu, 25 Mar 2021"
$fB10 = $wsHome.Range("B10").Font
$fB10.Name = "PingFang SC"
$fB10.Size = 12
$wsHome.Rows.Item(10).RowHeight = 19

# Row 11 (new) - Yourcontactdetails / 05 439 34875
$wsHome.Range("A11").Value = "Yourcontactdetails"
$fA11 = $wsHome.Range("A11").Font
$fA11.Name = "PingFang SC"
$fA11.Size = 12

$wsHome.Range("B11").Value = "05 439 34875"
$fB11 = $wsHome.Range("B11").Font
$fB11.Name = "PingFang SC"
$fB11.Size = 12
$wsHome.Rows.Item(11).RowHeight = 19

# Row 12 (new) - Installationdate / Fri, 26 Mar 2021, 8:00 am-12:00 pm
$wsHome.Range("A12").Value = "Installationdate"
$fA12 = $wsHome.Range("A12").Font
$fA12.Name = "PingFang SC"
$fA12.Size = 12

$wsHome.Range("B12").Value = "Fri, 26 Mar 2021, 8:00 am-12:00 pm"
$fB12 = $wsHome.Range("B12").Font
$fB12.Name = "PingFang SC"
$fB12.Size = 12
$wsHome.Rows.Item(12).RowHeight = 19

# Column B width widened to fit new longer content
$wsHome.Columns.Item(2).ColumnWidth = 41

# ---------------------------------------------------------------------------
# FeedBackPage (sheet6) - append addressOnMap / World Trade Center Mall row
# ---------------------------------------------------------------------------
$wsFeed.Range("A4").Value = "addressOnMap"
$wsFeed.Range("B4").Value = "World Trade Center Mall"
$wsFeed.Columns.Item(2).ColumnWidth = 23.75

# ---------------------------------------------------------------------------
# Selections per sheet (recorded cursor position at save time)
# ---------------------------------------------------------------------------
$wsSetUp.Range("D7").Select()
$wsLogin.Range("B4").Select()
$wsHelp.Range("B9").Select()
$wsFeed.Range("A4").Select()

# HomeMovePage becomes the active / visible tab, with B10 selected
$wsHome.Activate()
$wsHome.Range("B10").Select()

Write-Host "edit complete"
